$d = $word.ActiveDocument

# --- Part 1: tweak the two column widths on the "Content / Web Site" table ---
# Original grid: 3813 / 4106 (twips) -> target: 3823 / 4096 (twips)
# Column.Width is expressed in points (1 pt = 20 twips).
$t = $d.Tables.Item(1)
$t.Columns.Item(1).Width = 3823 / 20
$t.Columns.Item(2).Width = 4096 / 20

# --- Part 2: remove the trailing empty "BodyText" paragraph right after the table ---
$afterTable = $t.Range.End
$docEnd = $d.Content.End
$trailing = $d.Range($afterTable - 1, $docEnd)
$trailing.Delete()
